# Applies the "Update countries & provincias Spain" data refresh to the
# Pais worksheet: new totals for several countries, two countries whose
# case counts overtook their neighbours in the (sorted-by-total) table,
# and the refreshed "Datos actualizados" timestamp.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 23 de Abril de 2020 a las 02:22"

# Row 4: Estados Unidos
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 848717
$ws.Range("C4").Value = 29973
$ws.Range("D4").Value = 84048
$ws.Range("E4").Value = 717010
$ws.Range("F4").Value = 14016
$ws.Range("G4").Value = 2341
$ws.Range("H4").Value = 47659

# Row 8: Alemania
$ws.Range("A8").Value = "Alemania"
$ws.Range("B8").Value = 150648
$ws.Range("C8").Value = 2195
$ws.Range("D8").Value = 99400
$ws.Range("E8").Value = 45933
$ws.Range("F8").Value = 2908
$ws.Range("G8").Value = 229
$ws.Range("H8").Value = 5315

# Row 94: Crucero
$ws.Range("A94").Value = "Crucero"
$ws.Range("B94").Value = 712
$ws.Range("C94").Value = 0
$ws.Range("D94").Value = 645
$ws.Range("E94").Value = 54
$ws.Range("F94").Value = 4
$ws.Range("G94").Value = 0
$ws.Range("H94").Value = 13

# Row 120: Venezuela
$ws.Range("A120").Value = "Venezuela"
$ws.Range("B120").Value = 298
$ws.Range("C120").Value = 10
$ws.Range("D120").Value = 122
$ws.Range("E120").Value = 166
$ws.Range("F120").Value = 4
$ws.Range("G120").Value = 0
$ws.Range("H120").Value = 10

# Row 121: Mali
$ws.Range("A121").Value = "Mali"
$ws.Range("B121").Value = 293
$ws.Range("C121").Value = 35
$ws.Range("D121").Value = 73
$ws.Range("E121").Value = 203
$ws.Range("F121").Value = 0
$ws.Range("G121").Value = 3
$ws.Range("H121").Value = 17

# Row 125: Jamaica
$ws.Range("A125").Value = "Jamaica"
$ws.Range("B125").Value = 252
$ws.Range("C125").Value = 24
$ws.Range("D125").Value = 27
$ws.Range("E125").Value = 219
$ws.Range("F125").Value = 0
$ws.Range("G125").Value = 0
$ws.Range("H125").Value = 6

# Row 126: El Salvador
$ws.Range("A126").Value = "El Salvador"
$ws.Range("B126").Value = 237
$ws.Range("C126").Value = 12
$ws.Range("D126").Value = 63
$ws.Range("E126").Value = 167
$ws.Range("F126").Value = 3
$ws.Range("G126").Value = 0
$ws.Range("H126").Value = 7

# Row 142: Guayana Francesa
$ws.Range("A142").Value = "Guayana Francesa"
$ws.Range("B142").Value = 107
$ws.Range("C142").Value = 10
$ws.Range("D142").Value = 84
$ws.Range("E142").Value = 22
$ws.Range("F142").Value = 1
$ws.Range("G142").Value = 0
$ws.Range("H142").Value = 1

# Row 143: Liberia
$ws.Range("A143").Value = "Liberia"
$ws.Range("B143").Value = 101
$ws.Range("C143").Value = 0
$ws.Range("D143").Value = 20
$ws.Range("E143").Value = 73
$ws.Range("F143").Value = 0
$ws.Range("G143").Value = 0
$ws.Range("H143").Value = 8

# Row 144: Aruba
$ws.Range("A144").Value = "Aruba"
$ws.Range("B144").Value = 100
$ws.Range("C144").Value = 3
$ws.Range("D144").Value = 68
$ws.Range("E144").Value = 30
$ws.Range("F144").Value = 4
$ws.Range("G144").Value = 0
$ws.Range("H144").Value = 2

# Row 145: Bermudas
$ws.Range("A145").Value = "Bermudas"
$ws.Range("B145").Value = 99
$ws.Range("C145").Value = 1
$ws.Range("D145").Value = 39
$ws.Range("E145").Value = 55
$ws.Range("F145").Value = 10
$ws.Range("G145").Value = 0
$ws.Range("H145").Value = 5

# Row 151: Barbados
$ws.Range("A151").Value = "Barbados"
$ws.Range("B151").Value = 76
$ws.Range("C151").Value = 1
$ws.Range("D151").Value = 27
$ws.Range("E151").Value = 43
$ws.Range("F151").Value = 4
$ws.Range("G151").Value = 1
$ws.Range("H151").Value = 6

# Row 161: Libia
$ws.Range("A161").Value = "Libia"
$ws.Range("B161").Value = 60
$ws.Range("C161").Value = 1
$ws.Range("D161").Value = 15
$ws.Range("E161").Value = 44
$ws.Range("F161").Value = 0
$ws.Range("G161").Value = 0
$ws.Range("H161").Value = 1

# Row 206: Butan
$ws.Range("A206").Value = "Butan"
$ws.Range("B206").Value = 7
$ws.Range("C206").Value = 1
$ws.Range("D206").Value = 2
$ws.Range("E206").Value = 5
$ws.Range("F206").Value = 0
$ws.Range("G206").Value = 0
$ws.Range("H206").Value = 0

# Row 207: Mauritania
$ws.Range("A207").Value = "Mauritania"
$ws.Range("B207").Value = 7
$ws.Range("C207").Value = 0
$ws.Range("D207").Value = 6
$ws.Range("E207").Value = 0
$ws.Range("F207").Value = 0
$ws.Range("G207").Value = 0
$ws.Range("H207").Value = 1

# Row 208: Sahara Occidental
$ws.Range("A208").Value = "Sahara Occidental"
$ws.Range("B208").Value = 6
$ws.Range("C208").Value = 0
$ws.Range("D208").Value = 0
$ws.Range("E208").Value = 6
$ws.Range("F208").Value = 0
$ws.Range("G208").Value = 0
$ws.Range("H208").Value = 0
